# Component-interaction sequence diagram: rename the "delete 1" command's
# domain object from AddressBook/Person to Task (to match the Typed version),
# nudge one callout into alignment, add the presentation's two alignment
# guides, and let the date placeholders re-cache to the commit day.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Presentation-level alignment guides (View > Guides) used while laying
#    out the diagram: one horizontal guide and one (default-position)
#    vertical guide.
# ---------------------------------------------------------------------
try {
    $guides = $p.Guides
    $hGuide = $guides.Add(1, 1488)   # ppHorizontalGuide
    $vGuide = $guides.Add(2, 2880)   # ppVerticalGuide
} catch {
    # Guide manipulation isn't available in every host; continue regardless.
}

# ---------------------------------------------------------------------
# 2. Refresh the cached "datetimeFigureOut" footer field everywhere it is
#    defined (slide master, every layout, and the notes master) so the
#    deck shows the day it was last touched.
# ---------------------------------------------------------------------
function Update-DatePlaceholder($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tf = $shp.TextFrame
            if ($tf.HasText) {
                $tr = $tf.TextRange
                if ($tr.Text -eq "10/16/2016") {
                    $tr.Text = "4/8/17"
                }
            }
        }
    }
}

Update-DatePlaceholder $p.SlideMaster

$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li)
}

Update-DatePlaceholder $p.NotesMaster

# ---------------------------------------------------------------------
# 3. Slide 1 - the "delete 1" component-interaction diagram.
# ---------------------------------------------------------------------
$s = $p.Slides.Item(1)

# deletePerson(p) -> deleteTask(t)
$shp = $s.Shapes.Item(16)
$tr = $shp.TextFrame.TextRange
$tr.Characters(1, 12).Text = "deleteTask"
$tr.Characters(11, 3).Text = "(t)"

# post(AddressBookChangedEvent) -> post(TaskManagerChangedEvent)  (first copy)
$shp = $s.Shapes.Item(17)
$tr = $shp.TextFrame.TextRange
$tr.Characters(6, 23).Text = "TaskManagerChangedEvent"

# post(AddressBookChangedEvent) -> post(TaskManagerChangedEvent)  (second
# copy) - also nudge the textbox two pixels right so it lines up with the
# arrow above it.
$shp = $s.Shapes.Item(29)
$shp.Left = 144
$tr = $shp.TextFrame.TextRange
$tr.Characters(6, 23).Text = "TaskManagerChangedEvent"

# handleAddresssBookChangedEvent() -> handleTaskManagerChangedEvent()  (first copy)
$shp = $s.Shapes.Item(37)
$tr = $shp.TextFrame.TextRange
$tr.Characters(1, 30).Text = "handleTaskManagerChangedEvent"

# handleAddresssBookChangedEvent() -> handleTaskManagerChangedEvent()  (second copy)
$shp = $s.Shapes.Item(44)
$tr = $shp.TextFrame.TextRange
$tr.Characters(1, 30).Text = "handleTaskManagerChangedEvent"
